# Jesse's Week 10 logs — fill in the Task Summary Sheet and the
# Activity Log Summary Sheet with this week's data.

$wb = $excel.ActiveWorkbook

# --- TASK SUMMARY SHEET -----------------------------------------------
$ws3 = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Header: name + week number
$ws3.Range("C1").Value = "Jesse Hare"
$ws3.Range("E1").Value = 10

# Row 3: first task this week
$ws3.Range("A3").Value = "Project Build"
$ws3.Range("B3").Value = "Finalising feature set before next meeting"
$ws3.Range("C3").Value = 15
$ws3.Range("D3").Value = 15
$ws3.Range("E3").Value = 0

# Row 4: second task this week
$ws3.Range("A4").Value = "Project Build"
$ws3.Range("B4").Value = "Commence documentation writeup for searcher program"
$ws3.Range("C4").Value = 5
$ws3.Range("D4").Value = 5
$ws3.Range("E4").Value = 0

# --- ACTIVITY LOG SUMMARY SHEET ----------------------------------------
$ws4 = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")

# Name, right next to the sheet title
$ws4.Range("D1").Value = "Jesse Hare"

# Row 4: totals for the stage worked on this week
$ws4.Range("A4").Value = "Project Build"
$ws4.Range("B4").Value = 14
$ws4.Range("C4").Value = 6
